# fix: ajuste baseCiclcio e create generateFichas
# Append three new "Item" rows (A38:A40) to Sheet1, growing the used
# range from A1:A37 to A1:A40. Row 39 and 40 intentionally share the
# same item code ("0387-0199-03"), matching the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "0387-0199-01"
$ws.Range("A39").Value = "0387-0199-03"
$ws.Range("A40").Value = "0387-0199-03"

# Leave the selection on the last entered cell, as in the authored edit.
$ws.Range("A40").Select()
$excel.ActiveWindow.ScrollRow = 22
